$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Update the herbie_script text for the last task row (row 4, column E)
$ws.Range("E4").Value = 'verify url equals "https://mieweb.github.io/herbie/playgrounds/login.html"'

# Row 4 grows taller once the longer, wrapped text needs two lines
$ws.Rows.Item(4).RowHeight = 47.25

# Maximize-style window geometry change recorded by Excel on save
$excel.Width = 29040
$excel.Height = 15720
$wb.Windows.Item(1).WindowState = -4137
